$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new EUR->ARS quote row.
# Column A holds a date-looking string ("2025-09-29"); Excel's normal type
# inference would silently convert that to a date serial number when
# assigned through .Value. Force the cell to Text first so it stays a
# literal string (matching the workbook's existing inlineStr cells), then
# restore the cell style to "Normal" so no stray number-format styling is
# left behind on the cell.
$a47 = $ws.Range("A47")
$a47.NumberFormat = "@"
$a47.Value = "2025-09-29"
$a47.Style = "Normal"

$ws.Range("B47").Value = "15:20:20"
$ws.Range("C47").Value = "1.00 EUR = 1,625.3237"
